$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting rows 10-12 down to 11-13
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with data (a new weekly observation)
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44476
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100108
$ws.Range("H10").Value = "Tropicales y subtropicales"
$ws.Range("I10").Value = 100108001
$ws.Range("J10").Value = "Guayaba"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 1200
$ws.Range("O10").Value = 1200
$ws.Range("P10").Value = 1200
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región de Arica y Parinacota"
$ws.Range("S10").Value = 1200
$ws.Range("T10").Value = 1

$wb.Save()
